$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 239 (last existing data row) holds the template row (values + formatting)
# that gets duplicated into three new rows (240-242). Only column A (the date
# serial number) increments by 1 day for each new row; all other columns keep
# the same values as row 239.

$lastRow = 239
$newRows = 3
$template = $ws.Range("A$lastRow" + ":J$lastRow")

for ($i = 1; $i -le $newRows; $i++) {
    $r = $lastRow + $i
    $dst = $ws.Range("A$r" + ":J$r")

    # Copy the whole row (values + number formats/styles) from the template row.
    $template.Copy($dst)

    # Column A is the date; bump it by $i days relative to the template row.
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($lastRow, 1).Value2 + $i
}
